$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update report generation timestamp
$ws.Range("D5").Value = "Report Generated On: 08/26/2025 10:01 AM"

# Update total billed amount
$ws.Range("C8").Value = 412.36

# Clear Scope ID value (keep cell but remove text)
$ws.Range("G10").Value = ""

# Update pricing for line items
$ws.Range("H16").Value = 350.53
$ws.Range("H17").Value = 61.83
$ws.Range("H18").Value = 412.36
